$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume figures (and the Maker/TheSandbox/PaxDollar row rotation)
# for rows 2-51 per the latest GitHub Actions scrape.
# Values are written with a leading apostrophe so Excel stores them as literal text
# (matching the source inlineStr cells) instead of auto-coercing numeric-looking
# strings like "1.000" into the number 1. Style is reset to Normal afterwards so the
# quote-prefix flag Excel applies doesn't leave a stray cell style behind.

$ws.Range('D2').Value = '''30.367.61'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -2.08%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.905.19'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -2.73%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''0.9997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.02%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''238.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -2.64%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''0.9991'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -0.05%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.4759'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -2.24%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  -3.03%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.06695'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -3.49%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = '''  -3.67%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''103.08'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -4.89%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.07711'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -1.22%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''1.916.08'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -2.25%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = '''  -5.35%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''0.6758'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -3.83%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''260.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -7.49%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''30.393.94'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -2.04%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''0.9991'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -0.01%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''0.000007494'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -3.74%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''12.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -4.63%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''5.410'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -2.35%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''1.000'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  +0.04%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''6.316'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -3.40%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''9.448'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -4.45%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''164.55'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -2.53%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''19.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -4.89%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''2.056'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -6.44%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D29').Value = '''1.373'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -1.00%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''4.662'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -0.07%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''1.511'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -3.76%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''4.252'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -5.16%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''0.04777'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -3.24%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''0.7300'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -3.81%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''1.114'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -5.03%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''0.9984'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -0.14%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = '''  -0.78%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''0.01925'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -4.70%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''2.600'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -3.82%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''6.236'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -5.51%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''74.97'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -4.08%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''1.996'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  -6.30%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.8666'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -4.05%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''106.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -3.08%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = '''Maker'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = '''1.054.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +4.00%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = '''TheSandbox'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = '''https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = '''0.4254'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -4.71%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = '''PaxDollar'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = '''https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = '''0.9982'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -0.17%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''7.483'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -8.37%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''35.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -2.47%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''0.1198'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -4.82%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''8.920'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -4.87%  '
$ws.Range('E51').Style = 'Normal'
